$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "BA"
$ws.Range("C3").Value = "SMCI"
$ws.Range("C4").Value = "SHOP.TO"
$ws.Range("C5").Value = "BTC-USD"
$ws.Range("C6").Value = "I like cats"

$ws.Range("C7").Select()
